# Daily update at 8 AM UTC
# Appends the next day's row of data to the "Wins Over Time" sheet and
# moves the "latest row" date formatting down onto the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Previously-last row (20) reverts from the "date only, last row" format
# to the regular "date + time" format used by every other data row.
$ws.Range("A20").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New data row for the next day.
$ws.Range("A21").Value = 45970
$ws.Range("B21").Value = 45
$ws.Range("C21").Value = 51
$ws.Range("D21").Value = 51

# The new last row gets the "date only" formatting that row 20 used to have.
$ws.Range("A21").NumberFormat = "YYYY-MM-DD"
